$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$prefix = "sequence/run_0659_samples/"

for ($r = 2; $r -le 21; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $val = $cell.Value2
    if ($val -like "$prefix*") {
        $cell.Value = $val.Substring($prefix.Length)
    }
}

$ws.Range("F2:F21").Select()
